# Regenerate the "K" (strikeouts) column (column G) of the save-data sheet.
# The previous export wrote Strike# (TB) into column G by mistake; this
# recalculates the true K value for every game row and writes it back,
# along with refreshed std/mean style values already baked into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number (1-based, matching the worksheet) -> corrected K value.
$kValues = @{
    2 = 2
    3 = 1
    4 = 1
    5 = 1
    6 = 0
    7 = 0
    8 = 0
    9 = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 2
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 1
    20 = 3
    21 = 0
    22 = 2
    23 = 2
    24 = 0
    25 = 0
    26 = 0
    27 = 1
    28 = 3
    29 = 2
    30 = 1
    31 = 0
    32 = 3
    33 = 2
    35 = 1
    36 = 1
    37 = 2
    38 = 1
    39 = 1
    40 = 0
    41 = 1
    42 = 1
    43 = 5
    44 = 2
    45 = 3
    46 = 0
    48 = 1
    49 = 1
    50 = 2
    51 = 2
    52 = 2
    53 = 0
    54 = 0
    55 = 1
    56 = 1
    57 = 2
    59 = 1
    60 = 0
    61 = 0
    62 = 2
    63 = 1
    64 = 0
    65 = 1
    66 = 0
    67 = 2
    68 = 2
    69 = 0
    70 = 0
    71 = 2
    72 = 2
    73 = 1
    74 = 2
    75 = 1
    76 = 1
    77 = 2
    78 = 2
    79 = 2
    80 = 0
    81 = 3
    82 = 1
    83 = 1
    84 = 1
    85 = 3
    86 = 2
    87 = 2
    88 = 2
    89 = 1
    90 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
